# Update several mass/abundance estimates in the regression-results table.
# Each edit targets one specific table cell (by row/column) so that the
# correct occurrence is changed even though some numbers (e.g. "0.037",
# "0.038", "4.8 x 10", "1.9 x 10" ...) repeat elsewhere in the table.
#
# Note: Find/Replace executed directly on a Table.Cell's Range can end up
# acting on the first matching text anywhere in the document in this
# runtime, so we rebuild an equivalent Range via $d.Range(start, end)
# before running Find - that correctly scopes the search/replace to just
# the target cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-InCell {
    param($Doc, $Table, [int]$Row, [int]$Col, [string]$OldText, [string]$NewText)

    $cell = $Table.Cell($Row, $Col)
    $scoped = $Doc.Range($cell.Range.Start, $cell.Range.End)
    $found = $scoped.Find.Execute($OldText, $true, $false, $false, $false, $false, `
                                   $true, 0, $false, $NewText, 1)
    if (-not $found) {
        Write-Host "WARNING: text not found for row $Row col $Col : $OldText"
    }
}

# Diversity (Shannon H) - Intercept row
Replace-InCell $d $t 14 4 "7.9" "8.1"
Replace-InCell $d $t 14 5 "2.12" "2.09"
Replace-InCell $d $t 14 6 "0.037" "0.040"

# Diversity (Shannon H) - log(Size Class) row
Replace-InCell $d $t 15 4 "4.8" "4.9"
Replace-InCell $d $t 15 5 "4.25" "4.08"

# Diversity (Shannon H) - log(Size Class)^2 row
Replace-InCell $d $t 16 4 "8.9" "9.1"
Replace-InCell $d $t 16 5 "-4.15" "-4.03"

# Diversity (Shannon H) - Latitude row
Replace-InCell $d $t 17 4 "4.1" "4.2"
Replace-InCell $d $t 17 5 "-2.11" "-2.08"
Replace-InCell $d $t 17 6 "0.038" "0.041"

# Diversity (Shannon H) - Latitude^2 row
Replace-InCell $d $t 18 4 "5.3" "5.5"
Replace-InCell $d $t 18 5 "2.11" "2.08"
Replace-InCell $d $t 18 6 "0.038" "0.041"

# Diversity (Shannon H) - Depth row
Replace-InCell $d $t 19 3 "1.9" "2.1"
Replace-InCell $d $t 19 5 "0.99" "1.09"
Replace-InCell $d $t 19 6 "0.325" "0.281"

# Evenness (Pielou J) - log(Size Class) row
Replace-InCell $d $t 21 6 "0.053" "0.052"

# Evenness (Pielou J) - log(Size Class)^2 row
Replace-InCell $d $t 22 6 "0.069" "0.068"

# Evenness (Pielou J) - Depth row
Replace-InCell $d $t 25 6 "0.720" "0.723"
